$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Clear all existing hyperlinks on the sheet so row-shifted cells do not retain stale link targets
$ws.Hyperlinks.Delete()

# Column B got one character wider (51 -> 52 chars)
$ws.Columns.Item(2).ColumnWidth = 51.1

# Full row data for rows 2..14 (A..H), reflecting the new scrape snapshot at 2025-11-30 18:24:03
$rowsData = @(
  @("2025-11-30 18:24:03", "【AI系勉強会】「Google Gravity」開発事例発表者募集!個人開発をプレゼンしませんか?", "システム開発", "1,000 ~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5443957", 360, "🔥AI,Ai ◆開発"),
  @("2025-11-30 18:24:03", "【急募】BlockChainとSolidityに精通したプログラマー募集", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5443998", 303, "🔥AI,Ai"),
  @("2025-11-30 18:24:03", "顧客のBtoB向けの管理画面サービスの構築とAPI連携", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5444251", 220, "🔥API ◇管理"),
  @("2025-11-30 18:24:03", "【自動化】ニュースサイト情報をX(旧Twitter)へ投稿するシステム開発", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5444198", 205, "◆開発,システム開発 ◇サイト"),
  @("2025-11-30 18:24:03", "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5251319", 135, "◆ツール,スクレイピング ◇サイト"),
  @("2025-11-30 18:24:03", "【外国人大歓迎】【急募】ECツールの保守・バグ修正・機能追加エンジニア募集", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5444141", 68, "◆ツール"),
  @("2025-11-30 18:24:03", "【急募】革新的ペット向けECプラットフォーム開発エンジニア募集", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5443928", 68, "◆開発"),
  @("2025-11-30 18:24:03", "【急募】魅力的なWEBサイト制作のフリーランスを探しています!", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5444036", 45, "◇サイト"),
  @("2025-11-30 18:24:03", "wordpressレンダリングを妨げるリソースの除外", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5016989", 33, "○WordPress"),
  @("2025-11-30 18:24:03", "Salesforce Agentforceの構築・導入支援", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437485", 18, $null),
  @("2025-11-30 18:24:03", "comfyui(paperspace)でエロ動画のループ物を作成したいです。その方法を教えてください", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5444370", 10, $null),
  @("2025-11-30 18:24:03", "【急募】ミニPCでクラウドストレージ(nextcloud)とOpenWrtルータ化の依頼", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5444262", 10, $null),
  @("2025-11-30 18:24:03", "空き室情報を拾ってくスクリプト作成", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5444064", 10, $null)
)

$r = 2
foreach ($row in $rowsData) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]
  $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[5])
  $r = $r + 1
}
